$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2 through 28
# from serial date 45510 (2024-08-06) to 45511 (2024-08-07)
for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45510) {
        $cell.Value2 = 45511
    }
}
